$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_G"

$ws.Range("A2").Value = 70.652173913043484
$ws.Range("A3").Value = 70.380434782608688
$ws.Range("A4").Value = 70.380434782608688
$ws.Range("A5").Value = 70.108695652173907
$ws.Range("A6").Value = 69.83695652173914
$ws.Range("A7").Value = 70.108695652173907
$ws.Range("A8").Value = 71.739130434782609
$ws.Range("A9").Value = 72.282608695652172
$ws.Range("A10").Value = 72.554347826086953
$ws.Range("A11").Value = 72.282608695652172
$ws.Range("A12").Value = 70.652173913043484
$ws.Range("A13").Value = 71.739130434782609
$ws.Range("A14").Value = 72.554347826086953
$ws.Range("A15").Value = 72.554347826086953
$ws.Range("A16").Value = 72.554347826086953
$ws.Range("A17").Value = 71.467391304347828
$ws.Range("A18").Value = 71.739130434782609
$ws.Range("A19").Value = 71.739130434782609
$ws.Range("A20").Value = 71.739130434782609
$ws.Range("A21").Value = 72.282608695652172
$ws.Range("A22").Value = 72.554347826086953
$ws.Range("A23").Value = 69.021739130434781
$ws.Range("A24").Value = 70.380434782608688
$ws.Range("A25").Value = 69.565217391304344
$ws.Range("A26").Value = 71.195652173913047
$ws.Range("A27").Value = 71.739130434782609
$ws.Range("A28").Value = 72.282608695652172
$ws.Range("A29").Value = 73.369565217391312
$ws.Range("A30").Value = 73.097826086956516
$ws.Range("A31").Value = 73.097826086956516
$ws.Range("A32").Value = 69.565217391304344
$ws.Range("A33").Value = 69.83695652173914
$ws.Range("A34").Value = 70.923913043478265
$ws.Range("A35").Value = 70.652173913043484
$ws.Range("A36").Value = 70.108695652173907
$ws.Range("A37").Value = 71.195652173913047
$ws.Range("A38").Value = 69.021739130434781
$ws.Range("A39").Value = 69.021739130434781
$ws.Range("A40").Value = 69.565217391304344
$ws.Range("A41").Value = 72.282608695652172
$ws.Range("A42").Value = 72.282608695652172
$ws.Range("A43").Value = 72.282608695652172
$ws.Range("A44").Value = 71.739130434782609
$ws.Range("A45").Value = 69.565217391304344
$ws.Range("A46").Value = 70.380434782608688
$ws.Range("A47").Value = 69.565217391304344
$ws.Range("A48").Value = 70.380434782608688
$ws.Range("A49").Value = 69.83695652173914
